# Apply the updated cryptocurrency market data (price + 1h volume change) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.635.18'

$ws.Range('D3').Value = '1.563.52'
$ws.Range('E3').Value = '  +0.14%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').Value = '''210.10'
$ws.Range('E5').Value = '  -0.13%  '

$ws.Range('D6').Value = '''0.510'
$ws.Range('E6').Value = '  +3.90%  '

$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.23%  '

$ws.Range('D8').Value = '''24.89'
$ws.Range('E8').Value = '  +5.29%  '

$ws.Range('E9').Value = '  +0.97%  '

$ws.Range('E10').Value = '  -0.06%  '

$ws.Range('D11').Value = '''0.0896'
$ws.Range('E11').Value = '  +0.30%  '

$ws.Range('D12').Value = '1.787.74'
$ws.Range('E12').Value = '  +0.17%  '

$ws.Range('D13').Value = '1.565.33'
$ws.Range('E13').Value = '  +0.08%  '

$ws.Range('D14').Value = '28.682.82'
$ws.Range('E14').Value = '  +1.36%  '

$ws.Range('E15').Value = '  +0.94%  '

$ws.Range('E16').Value = '  -0.33%  '

$ws.Range('D17').Value = '''61.46'
$ws.Range('E17').Value = '  +0.82%  '

$ws.Range('D18').Value = '''227.50'
$ws.Range('E18').Value = '  -0.17%  '

$ws.Range('E19').Value = '  -0.51%  '

$ws.Range('D20').Value = '0.0₃0679'
$ws.Range('E20').Value = '  +0.24%  '

$ws.Range('E21').Value = '  -0.22%  '

$ws.Range('E22').Value = '  -0.06%  '

$ws.Range('D23').Value = '''9.01'
$ws.Range('E23').Value = '  +1.49%  '

$ws.Range('E24').Value = '  +1.47%  '

$ws.Range('D25').Value = '''151.58'
$ws.Range('E25').Value = '  +0.94%  '

$ws.Range('D27').Value = '''14.76'
$ws.Range('E27').Value = '  -0.56%  '

$ws.Range('E28').Value = '  -0.14%  '

$ws.Range('E29').Value = '  -1.59%  '

$ws.Range('D30').Value = '''0.0457'
$ws.Range('E30').Value = '  -3.99%  '

$ws.Range('E31').Value = '  -1.62%  '

$ws.Range('E32').Value = '  +0.26%  '

$ws.Range('D33').Value = '1.398.99'
$ws.Range('E33').Value = '  +1.50%  '

$ws.Range('D34').Value = '''2.98'
$ws.Range('E34').Value = '  -3.01%  '

$ws.Range('D35').Value = '''1.03'
$ws.Range('E35').Value = '  -3.07%  '

$ws.Range('E36').Value = '  -1.83%  '

$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '''2.66'
$ws.Range('E37').Value = '  +1.20%  '

$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '''2.30'
$ws.Range('E38').Value = '  -2.04%  '

$ws.Range('E39').Value = '  -0.47%  '

$ws.Range('D40').Value = '''1.94'
$ws.Range('E40').Value = '  +0.76%  '

$ws.Range('D41').Value = '''0.517'
$ws.Range('E41').Value = '  -0.55%  '

$ws.Range('E42').Value = '  -0.22%  '

$ws.Range('D43').Value = '''0.767'
$ws.Range('E43').Value = '  -1.57%  '

$ws.Range('D44').Value = '''0.0461'
$ws.Range('E44').Value = '  -2.23%  '

$ws.Range('D45').Value = '''63.83'
$ws.Range('E45').Value = '  +2.93%  '

$ws.Range('E46').Value = '  -1.85%  '

$ws.Range('D47').Value = '1.699.73'
$ws.Range('E47').Value = '  +0.16%  '

$ws.Range('D48').Value = '''0.863'
$ws.Range('E48').Value = '  -6.00%  '

$ws.Range('D49').Value = '''84.73'
$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('D50').Value = '''42.40'
$ws.Range('E50').Value = '  +5.96%  '

$ws.Range('E51').Value = '  -0.42%  '
